$d = $word.ActiveDocument
$nbsp = [char]0x00a0

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Host "WARNING: not found -> $find"
    }
}

# "FAST" heading -> "Скорость"
Replace-Text "FAST" "Скорость"

# InstantPay teaser paragraph
Replace-Text ("We work daily to expand the offerings of SmartCash, from new mobile wallets to exchanges and community frameworks." + $nbsp + "InstantPay is just the tip of the iceberg.") "Каждый день мы работаем над тем, чтобы сделать SmartCash лучше – добавляем новые функции и методы оптимизации, разрабатываем мобильные приложения и сотрудничаем с биржами. InstantPay – это только вершина айсберга."

# "INCLUSIVE" heading -> "Вовлеченность"
Replace-Text "INCLUSIVE" "Вовлеченность"

# Inclusive paragraph
Replace-Text ("The entire foundation of SmartCash is based upon community adoption, growth and inclusion." + $nbsp + "Everyone is a SmartCash team member.") "Весь мир SmartCash выстраивается вокруг сообщества, которое принимает решения и помогает в развитии монеты. Каждый может стать частью команды."

# "JOIN THE SMARTCASH COMMUNITY" heading
Replace-Text "JOIN THE SMARTCASH COMMUNITY" "Присоединяйтесь к сообществу SmartCash"

# Community intro paragraph
Replace-Text "We have many different places and methods to discuss, learn and flourish together. Below are a few options." "Присоединяйтесь к обсуждению на различных платформах. Ниже представлены платформы, где вы можете узнать о наших новостях или обсудить важные темы:"

# " SERVICES" heading (leading nbsp)
Replace-Text ($nbsp + "SERVICES") "Сервисы"

# Services disclaimer paragraph
Replace-Text "DISCLAIMER: This list is provided for informational purposes only. Services listed here have not been evaluated or endorsed by the SmartCash developers and no guarantees are made as to the accuracy of this information. Please exercise discretion when using third-party services." "ОТКАЗ ОТ ОТВЕТСТВЕННОСТИ: Данный список носит ознакомительный характер. Сервисы, перечисленные здесь, не были оценены или проверены разработчиками SmartCash, поэтому гарантий относительно точности данной информации нет. Пожалуйста, будьте внимательны при использовании сторонних сервисов."

# " PROJECTS" heading (leading nbsp)
Replace-Text ($nbsp + "PROJECTS") "Проекты"

# " WALLETS" heading (leading nbsp)
Replace-Text ($nbsp + "WALLETS") "Кошельки"

# Fast wallet description
Replace-Text "This is a fast wallet that does not require the blockchain download. Wallet will not start SmartNodes, but a future release will add that feature." "Быстрый кошелёк, который не требует загрузки блокчейна. Функция запуска SmartNodes будет добавлена в будущих релизах."

# "NODE CLIENT" heading
Replace-Text "NODE CLIENT" "Node-клиент"

# Node client description (up to, not including, the hyperlinked "here")
Replace-Text ("This Node Client Supports SmartNode activation. To assist in speeding up the syncing of your Node Client read more" + $nbsp) "Node-клиент поддерживает запуск SmartNode. Для ускорения синхронизации вашего Node-клиента, пожалуйста, узнайте больше по "

# Hyperlink text "here" -> "ссылке." (handled separately to avoid touching the
# other, non-hyperlinked occurrence of the word "here" elsewhere in the doc,
# and to keep the hyperlink's run formatting as intact as possible)
$hyperlinkRange = $d.Hyperlinks(1).Range
$linkStart = $hyperlinkRange.Start
foreach ($ch in 1..40) {
    $probe = $d.Range($linkStart, $linkStart + 4)
    if ($probe.Text -eq "here") { break }
    $linkStart = $linkStart + 1
}
$linkRange = $d.Range($linkStart, $linkStart + 4)
$savedColor = $linkRange.Font.Color
$savedUnderline = $linkRange.Font.Underline
$savedBold = $linkRange.Font.Bold
$savedItalic = $linkRange.Font.Italic
$savedNameAscii = $linkRange.Font.NameAscii
$savedNameFarEast = $linkRange.Font.NameFarEast
$savedNameOther = $linkRange.Font.NameOther
$savedSize = $linkRange.Font.Size
$linkRange.Text = "ссылке."
$newLinkRange = $d.Range($linkStart, $linkStart + 7)
$newLinkRange.Font.NameAscii = $savedNameAscii
$newLinkRange.Font.NameFarEast = $savedNameFarEast
$newLinkRange.Font.NameOther = $savedNameOther
$newLinkRange.Font.Color = $savedColor
$newLinkRange.Font.Underline = $savedUnderline
$newLinkRange.Font.Bold = $savedBold
$newLinkRange.Font.Italic = $savedItalic
$newLinkRange.Font.Size = $savedSize

# " EXCHANGES" heading (leading nbsp)
Replace-Text ($nbsp + "EXCHANGES") "Биржи"

# " LOOKING FOR THE MAILING LIST?" heading (leading nbsp kept as a regular leading space)
Replace-Text ($nbsp + "LOOKING FOR THE MAILING LIST?") (" Хотите быть в курсе событий?")

# Mailing list intro paragraph
Replace-Text "Want to keep up with developments, news and updates from the SmartCash team?" "Желаете всегда быть в курсе событий, новостей и обновлений от команды SmartCash? Присоединяйтесь к нашей рассылке!"

# "Email Address " label (trailing nbsp)
Replace-Text ("Email Address" + $nbsp) "Ваш Email"

# "Subscribe" button
Replace-Text "Subscribe" "Подписаться"
